$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Row 39: Generation Friends / Saul Austerlitz
$ws.Range("A39").Value = "Generation Friends"
$ws.Range("B39").Value = "Saul Austerlitz"

$ws.Range("C2").Copy()
$ws.Range("C39:D40").PasteSpecial(-4122)  # xlPasteFormats, reuse the existing date style

$ws.Range("C39").Value = "3/17/2020"
$ws.Range("D39").Value = "3/20/2020"
$ws.Range("E39").Value = "friends;tv;hollywood;history;analysis"
$ws.Range("F39").Value = "Audio"
$ws.Range("G39").Value = "10 Hours 53 Mins"

# Row 40: American Icon / Bryce G. Hoffman
$ws.Range("A40").Value = "American Icon"
$ws.Range("B40").Value = "Bryce G. Hoffman"
$ws.Range("C40").Value = "3/8/2020"
$ws.Range("D40").Value = "3/21/2020"
$ws.Range("E40").Value = "alan mulally;ford;great recession;business;turn around; history;success"
$ws.Range("F40").Value = "Hard Copy"
$ws.Range("G40").Value = "398 Pages"

$excel.CutCopyMode = $false

$ws.Range("A41").Select()
